$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The "jobs" sheet stores every column except A (jobNumber, numeric)
# and L (_isDeleted, boolean) as literal text -- even values that look
# like numbers or dates (e.g. "460", "0", "8/31/2022"). Excel's COM
# Range.Value setter auto-detects numbers/dates, so every text column
# that is about to receive a value is first force-formatted as Text
# ("@") to keep those look-alike values stored as strings, matching
# the rest of the sheet.
# NOTE: each NumberFormat assignment below targets a single contiguous
# area (never a comma-separated multi-area range) on purpose.
# ------------------------------------------------------------------
$ws.Range("B27:C27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("I27:J27").NumberFormat = "@"
$ws.Range("O27").NumberFormat = "@"
$ws.Range("B28:O33").NumberFormat = "@"

# ---- Row 27: update the existing record ----
$ws.Range("B27").Value = "06D3136CC3600"
$ws.Range("C27").Value = "NO TAG"
$ws.Range("E27").Value = "0"
$ws.Range("G27").Value = "82523"
$ws.Range("I27").Value = ""
$ws.Range("J27").Value = "ravi"
$ws.Range("O27").Value = "NO"

# ---- Row 28: new record ----
$ws.Range("A28").Value = 71280
$ws.Range("B28").Value = "O6E7265 310"
$ws.Range("C28").Value = "6EABBC2L2H30"
$ws.Range("D28").Value = "MULTI"
$ws.Range("E28").Value = "?"
$ws.Range("F28").Value = "GOOD"
$ws.Range("G28").Value = "NA"
$ws.Range("H28").Value = "NO"
$ws.Range("I28").Value = ""
$ws.Range("J28").Value = "ravi"
$ws.Range("K28").Value = "8/31/2022"
$ws.Range("L28").Value = $false
$ws.Range("M28").Value = "N/A"
$ws.Range("N28").Value = "N/A"
$ws.Range("O28").Value = "NO"

# ---- Row 29: new record ----
$ws.Range("A29").Value = 71281
$ws.Range("B29").Value = "NRA0400TFD"
$ws.Range("C29").Value = "CT97E15052"
$ws.Range("D29").Value = "460"
$ws.Range("E29").Value = "0"
$ws.Range("F29").Value = "?"
$ws.Range("G29").Value = "N/A"
$ws.Range("H29").Value = "NO"
$ws.Range("I29").Value = ""
$ws.Range("J29").Value = "ravi"
$ws.Range("K29").Value = "9/2/2022"
$ws.Range("L29").Value = $false
$ws.Range("M29").Value = "N/A"
$ws.Range("N29").Value = "N/A"
$ws.Range("O29").Value = "NO"

# ---- Row 30: new record ----
$ws.Range("A30").Value = 71282
$ws.Range("B30").Value = "NRB 0400TFD"
$ws.Range("C30").Value = "CT0400TFD"
$ws.Range("D30").Value = "460"
$ws.Range("E30").Value = "0"
$ws.Range("F30").Value = "BAD"
$ws.Range("G30").Value = "N/A"
$ws.Range("H30").Value = "NO"
$ws.Range("I30").Value = ""
$ws.Range("J30").Value = "ravi"
$ws.Range("K30").Value = "9/2/2022"
$ws.Range("L30").Value = $false
$ws.Range("M30").Value = "N/A"
$ws.Range("N30").Value = "N/A"
$ws.Range("O30").Value = "NO"

# ---- Row 31: new record ----
$ws.Range("A31").Value = 71283
$ws.Range("B31").Value = "3DP3R12METFD"
$ws.Range("C31").Value = "14162424R"
$ws.Range("D31").Value = "460"
$ws.Range("E31").Value = "1"
$ws.Range("F31").Value = "GOOD"
$ws.Range("G31").Value = "N/A"
$ws.Range("H31").Value = "NO"
$ws.Range("I31").Value = ""
$ws.Range("J31").Value = "ravi"
$ws.Range("K31").Value = "9/2/2022"
$ws.Range("L31").Value = $false
$ws.Range("M31").Value = "N/A"
$ws.Range("N31").Value = "N/A"
$ws.Range("O31").Value = "NO"

# ---- Row 32: new record ----
$ws.Range("A32").Value = 71284
$ws.Range("B32").Value = "4RL2150ATSK"
$ws.Range("C32").Value = "21F60072R"
$ws.Range("D32").Value = "MULTI"
$ws.Range("E32").Value = "0"
$ws.Range("F32").Value = "GOOD"
$ws.Range("G32").Value = "N/A"
$ws.Range("H32").Value = "NO"
$ws.Range("I32").Value = ""
$ws.Range("J32").Value = "ravi"
$ws.Range("K32").Value = "9/2/2022"
$ws.Range("L32").Value = $false
$ws.Range("M32").Value = "N/A"
$ws.Range("N32").Value = "N/A"
$ws.Range("O32").Value = "NO"

# ---- Row 33: new record ----
$ws.Range("A33").Value = 71285
$ws.Range("B33").Value = "3DA3A060ETFC"
$ws.Range("C33").Value = "ET06G02970R"
$ws.Range("D33").Value = "208"
$ws.Range("E33").Value = "0"
$ws.Range("F33").Value = "BAD"
$ws.Range("G33").Value = "N/A"
$ws.Range("H33").Value = "NO"
$ws.Range("I33").Value = ""
$ws.Range("J33").Value = "ravi"
$ws.Range("K33").Value = "9/2/2022"
$ws.Range("L33").Value = $false
$ws.Range("M33").Value = "N/A"
$ws.Range("N33").Value = "N/A"
$ws.Range("O33").Value = "NO"
